# Update cryptos list figures (price / volume columns) and fix the
# ordering of the PaxDollar / InjectiveProtocol rows (48 & 49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new D (price) value. Column D holds plain text in the
# source workbook (e.g. "34.602.24", "224.85"); many of the new values
# would otherwise be auto-coerced to numbers by Excel's input parser, so
# each write is wrapped with a text NumberFormat that gets reset right
# afterwards (keeps the cell's style index unchanged).
$priceUpdates = @{
    2  = "34.674.06"
    3  = "1.790.59"
    5  = "224.78"
    6  = "0.562"
    8  = "32.64"
    10 = "0.0670"
    12 = "2.046.55"
    13 = "11.07"
    14 = "1.789.54"
    15 = "0.635"
    16 = "34.592.84"
    17 = "4.30"
    18 = "68.93"
    19 = "254.51"
    22 = "10.42"
    23 = "4.24"
    25 = "159.53"
    26 = "16.41"
    27 = "7.10"
    33 = "3.59"
    34 = "1.86"
    35 = "1.444.68"
    38 = "0.629"
    39 = "83.22"
    42 = "0.902"
    47 = "1.947.54"
    48 = "12.00"
    49 = "1.00"
    50 = "103.26"
}

# Map of row -> new E (volume) value. These already contain padding
# spaces and a trailing "%" so Excel's parser leaves them as text.
$volumeUpdates = @{
    2  = "  +2.25%  "
    3  = "  +0.33%  "
    4  = "  -0.07%  "
    5  = "  -0.07%  "
    6  = "  +1.25%  "
    7  = "  -0.02%  "
    8  = "  +6.14%  "
    9  = "  +2.23%  "
    10 = "  +1.45%  "
    11 = "  +1.35%  "
    12 = "  +0.23%  "
    13 = "  +10.78%  "
    14 = "  +0.29%  "
    15 = "  +0.86%  "
    16 = "  +2.02%  "
    17 = "  +2.58%  "
    18 = "  +0.58%  "
    19 = "  +1.22%  "
    20 = "  +3.65%  "
    21 = "  +0.05%  "
    22 = "  +1.19%  "
    23 = "  +0.67%  "
    24 = "  -1.40%  "
    25 = "  +0.71%  "
    26 = "  -0.43%  "
    27 = "  +2.34%  "
    28 = "  -0.02%  "
    29 = "  -0.09%  "
    30 = "  -0.84%  "
    31 = "  +0.69%  "
    32 = "  +0.21%  "
    33 = "  +0.97%  "
    34 = "  +3.64%  "
    35 = "  -2.72%  "
    36 = "  +0.12%  "
    37 = "  +2.48%  "
    38 = "  -0.43%  "
    39 = "  -0.17%  "
    40 = "  +4.44%  "
    41 = "  -0.08%  "
    42 = "  +1.60%  "
    43 = "  +0.30%  "
    44 = "  -0.87%  "
    45 = "  +2.40%  "
    46 = "  -1.62%  "
    47 = "  +0.42%  "
    48 = "  +1.29%  "
    49 = "  -0.03%  "
    50 = "  +5.59%  "
    51 = "  +6.22%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}

# Rows 48 and 49 swap coin name/link (new price/volume already set above).
$ws.Cells.Item(48, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"

$ws.Cells.Item(49, 2).Value = "PaxDollar"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
